{"js": "// Replace the placeholder ID text (this also removes the trailing space\n// run, since the replacement text + following space are merged into a\n// single new run by the text-replace operation).\nconst body = context.document.body;\nconst results = body.search(\"**ID__AFFARS_pgi_5307_topic_8__ID** \", {\n  matchCase: true,\n  matchWildcards: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"**ID__AFFARS_AFMC_PGI_5307__ID**\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// Update the first paragraph's indentation and add a paragraph border\n// (top/left/bottom/right, 5pt space from text, no line drawn).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.leftIndent = 11.25; // 225 twips (20 twips per point)\n\nconst borders = firstParagraph.borders;\nborders.load(\"items\");\nawait context.sync();\n\n// items order is Top, Left, Bottom, Right.\nborders.items[0]._omSet(\"DistanceFromTop\", 5, \"Border\");\nborders.items[1]._omSet(\"DistanceFromLeft\", 5, \"Border\");\nborders.items[2]._omSet(\"DistanceFromBottom\", 5, \"Border\");\nborders.items[3]._omSet(\"DistanceFromRight\", 5, \"Border\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the placeholder ID text. Matching through the trailing space and\n# replacing with the new ID (no trailing space) merges the two original\n# runs into a single new run, matching the target structure.\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"**ID__AFFARS_pgi_5307_topic_8__ID** \"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"**ID__AFFARS_AFMC_PGI_5307__ID**\"\n$find.Execute(\n  $find.Text,    # FindText\n  $false,        # MatchCase\n  $false,        # MatchWholeWord\n  $false,        # MatchWildcards\n  $false,        # MatchSoundsLike\n  $false,        # MatchAllWordForms\n  $true,         # Forward\n  $wdFindContinue,\n  $false,        # Format\n  $find.Replacement.Text,\n  $wdReplaceAll\n) | Out-Null\n\n# Update the first paragraph's indentation and add a paragraph border\n# (top/left/bottom/right, 5pt space from text, no line drawn).\n$p = $d.Paragraphs(1)\n$p.LeftIndent = 11.25  # 225 twips (20 twips per point)\n\n$borders = $p.Borders\n$borders.DistanceFromTop = 5\n$borders.DistanceFromLeft = 5\n$borders.DistanceFromBottom = 5\n$borders.DistanceFromRight = 5\n"}
